# Apply the "Changes in Dox and Box" edit to the Request sheet:
#   - Row 2 (F-000001161 / ACCOUNT MANUAL): Remark/Type Of Retrival/Type of
#     Selivery change from "dfs" / "Physical-Returnable" / "Urgent" to
#     "jhjh" / "Digital (Scan)" / "Standard".
#   - Row 3 (F-000001160 / ACCOUNT MANUAL / dfs / Physical-Returnable /
#     Urgent) is removed entirely, shrinking the used range to A1:E2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Request")

# Update the Remark / Type Of Retrival / Type of Selivery cells on row 2.
$ws.Range("C2").Value = "jhjh"
$ws.Range("D2").Value = "Digital (Scan)"
$ws.Range("E2").Value = "Standard"

# Delete row 3 entirely - remaining rows (none, here) shift up and the
# sheet's used range/dimension shrinks from A1:E3 to A1:E2.
$ws.Rows("3:3").Delete()
